$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data row (old row 3) down to make room for 5 new language rows.
# This also clones formatting (style) from the surrounding rows, matching how
# Excel/LibreOffice would propagate formatting on a row insert.
$ws.Rows("3:7").Insert()

# Copy formatting (styles) from row 2 into the newly inserted rows 3-7 so the
# new data rows pick up the same cell style indices as row 2 / the shifted row.
$ws.Range("A2:E2").Copy()
$ws.Range("A3:E7").PasteSpecial()

# Fill in the new language rows (French, Arabic, Kannada, Hindi, Tamil).
$ws.Range("A3").Value = "fra"
$ws.Range("B3").Value = "French"
$ws.Range("C3").Value = "Indo-European"
$ws.Range("D3").Value = "français"

$ws.Range("A4").Value = "ara"
$ws.Range("B4").Value = "Arabic"
$ws.Range("C4").Value = "الهندو أوروبية"
$ws.Range("D4").Value = "عربي"

$ws.Range("A5").Value = "kan"
$ws.Range("B5").Value = "ಕನ್ನಡ"
$ws.Range("C5").Value = "ಇಂಡೋ-ಯುರೋಪಿಯನ್"
$ws.Range("D5").Value = "Kannada"

$ws.Range("A6").Value = "hin"
$ws.Range("B6").Value = "हिन्दी"
$ws.Range("C6").Value = "भारोपीय"
$ws.Range("D6").Value = "Hindi"

$ws.Range("A7").Value = "tam"
$ws.Range("B7").Value = "தமிழ்"
$ws.Range("C7").Value = "இந்தோ-ஐரோப்பிய"
$ws.Range("D7").Value = "Tamil"

# "is_active" column (E) must stay a text "TRUE" (shared string), not a boolean,
# so copy it down from the existing TRUE cell instead of assigning Value directly.
for ($r = 3; $r -le 7; $r++) {
    $ws.Range("E2").Copy()
    $ws.Range("E$r").PasteSpecial()
}

# Leave the active selection on A9, matching the post-edit workbook state.
[void]$ws.Range("A9").Select()
